$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 3.286832544864788
$ws.Range("C2").Value = 1.655778082260271
$ws.Range("D2").Value = 0.1494219747398047
$ws.Range("E2").Value = 0.4942365360607697
$ws.Range("G2").Value = 5.586269137925634

$ws.Range("B3").Value = 0.6606524410359556
$ws.Range("C3").Value = 1.655778082260271
$ws.Range("D3").Value = 22.3905356188092
$ws.Range("E3").Value = 10.19245300693656
$ws.Range("G3").Value = 34.89941914904198
